# Auto-generated edit script: fills in previously-empty table cells
# with text content, matching the target revision of the document.
$d = $word.ActiveDocument

# Helper: fills an (already empty) table cell paragraph with one or more
# runs of text by repeatedly inserting after the collapsing range end.
# This preserves the existing <w:pPr> (alignment/style) of the paragraph.
function Set-CellRuns {
    param($table, $rowIdx, $cellIdx, $texts)
    $cell = $d.Tables.Item($table).Rows.Item($rowIdx).Cells.Item($cellIdx)
    $pr = $cell.Range.Paragraphs.Item(1).Range
    foreach ($txt in $texts) {
        $pr.InsertAfter($txt)
        $pr.Collapse(0)
    }
}

# Helper: fills an (already empty) table cell paragraph with several
# explicit <w:r> runs via InsertXML, preserving the paragraph identity
# (w14:paraId etc.) and its <w:pPr>. Used where the source document
# models the text as more than one run (e.g. a date typed, then a time
# appended in a second pass).
function Set-CellRunsXml {
    param($table, $rowIdx, $cellIdx, $paraId, $rsid, $pPrXml, $runsXml)
    $cell = $d.Tables.Item($table).Rows.Item($rowIdx).Cells.Item($cellIdx)
    $p = $cell.Range.Paragraphs.Item(1)
    $pr = $p.Range
    $pAttrs = 'w14:paraId="' + $paraId + '" w14:textId="77777777" w:rsidR="' + $rsid + '" w:rsidRDefault="' + $rsid + '" w:rsidP="00B5507D"'
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p ' + $pAttrs + '>' + $pPrXml + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pr.InsertXML($xml) | Out-Null
}

Set-CellRuns 2 6 4 @("0,67h")
Set-CellRuns 2 6 6 @("0,67h")
Set-CellRuns 2 13 8 @("66,66")
Set-CellRuns 2 14 4 @("0")
Set-CellRuns 2 14 6 @("0")
Set-CellRuns 2 14 8 @("0")
Set-CellRuns 2 15 4 @("1")
Set-CellRuns 2 15 6 @("1")
Set-CellRuns 2 15 8 @("33,34")
Set-CellRuns 2 16 4 @("3")
Set-CellRuns 2 16 6 @("3")
Set-CellRuns 2 21 8 @("66,66")
Set-CellRuns 2 22 4 @("0")
Set-CellRuns 2 22 6 @("0")
Set-CellRuns 2 22 8 @("0")
Set-CellRuns 2 23 4 @("1")
Set-CellRuns 2 23 6 @("1")
Set-CellRuns 2 23 8 @("33,34")
Set-CellRuns 2 24 4 @("3")
Set-CellRuns 2 24 6 @("3")
Set-CellRuns 3 5 1 @("1A")
Set-CellRuns 3 5 2 @("Test")
Set-CellRuns 3 5 3 @("10/Apr/20 21:30")
Set-CellRuns 3 5 4 @("-")
Set-CellRunsXml 3 5 5 '6843594F' '0042485E' '<w:pPr><w:pStyle w:val="FormText"/><w:snapToGrid w:val="0"/></w:pPr>' '<w:r><w:t>10/Apr/20</w:t></w:r><w:r><w:t xml:space="preserve"> 22:10</w:t></w:r>'
Set-CellRuns 3 5 6 @("0,67h")
Set-CellRuns 4 12 2 @("1A")
Set-CellRuns 4 12 4 @("10/Apr/20")
Set-CellRuns 4 12 6 @("3")
Set-CellRuns 4 12 8 @("Checking")
Set-CellRuns 4 12 10 @("1")
Set-CellRuns 4 12 12 @("1")
Set-CellRuns 4 12 14 @("0,08h")
Set-CellRuns 4 12 16 @("X")
Set-CellRuns 4 13 2 @("Service was not validating invalid arrays.")

Write-Host "Edit complete"